# Add 4 new Employer's Representative's Instruction document rows
# ("Removal of AtoNs" ERI#01-#04) to the "Documents" sheet, each with
# its own hyperlink, mirroring the commit "Add files via upload".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Documents")

$docType = "EMPLOYER'S REPRESENTATIVE'S INSTRUCTION"
$contractId = "DMCA-01"
$docTitle = "Removal of AtoNs"
$url = "https://dubaiholding-my.sharepoint.com/:b:/r/personal/arun_naidu_dhre_ae/Documents/Documents/NME/MARINE/DMS%20148763%20-%20Marine%20Works%20Including%20Dredging%20and%20Land%20Reclamation%20(Jan%20De%20Nul)/EMPLOYERS%20REPRESENTATIVES%20INSTRUCTION/ERI%2301/NKL-LT-12161%20-%20ERI%2301%20-%20Additional%20Sand%20Stockpiles.pdf?csf=1&web=1&e=7f8NQc"

# --- Row 4 : ERI1.0 (filled completely, left to right) ---
$ws.Cells.Item(4,1).Value = "ERI1.0"
$ws.Cells.Item(4,2).Value = $contractId
$ws.Cells.Item(4,3).Value = $docType
$ws.Cells.Item(4,4).Value = "ERI # 01.0"
$ws.Cells.Item(4,5).Value = $docTitle
$ws.Cells.Item(4,6).Value = "NKL-LT-12161 - ERI#01 - Removal of AtoNs.pdf"

# --- Rows 5-7 : columns A-E first ---
$ws.Cells.Item(5,1).Value = "ERI1.1"
$ws.Cells.Item(5,2).Value = $contractId
$ws.Cells.Item(5,3).Value = $docType
$ws.Cells.Item(5,4).Value = "ERI # 01.1"
$ws.Cells.Item(5,5).Value = $docTitle

$ws.Cells.Item(6,1).Value = "ERI1.2"
$ws.Cells.Item(6,2).Value = $contractId
$ws.Cells.Item(6,3).Value = $docType
$ws.Cells.Item(6,4).Value = "ERI # 01.2"
$ws.Cells.Item(6,5).Value = $docTitle

$ws.Cells.Item(7,1).Value = "ERI1.3"
$ws.Cells.Item(7,2).Value = $contractId
$ws.Cells.Item(7,3).Value = $docType
$ws.Cells.Item(7,4).Value = "ERI # 01.3"
$ws.Cells.Item(7,5).Value = $docTitle

# --- Column F (file name) for rows 5-7 ---
$ws.Cells.Item(5,6).Value = "NKL-LT-12161 - ERI#02 - Removal of AtoNs.pdf"
$ws.Cells.Item(6,6).Value = "NKL-LT-12161 - ERI#03 - Removal of AtoNs.pdf"
$ws.Cells.Item(7,6).Value = "NKL-LT-12161 - ERI#04 - Removal of AtoNs.pdf"

# --- Column G (hyperlink), H (upload date), I (version) for rows 4-7 ---
$dates = @("2024-05-12", "2024-05-13", "2024-05-14", "2024-05-15")
for ($r = 4; $r -le 7; $r++) {
    $cell = $ws.Cells.Item($r,7)
    $cell.Value = $url
    $ws.Hyperlinks.Add($cell, $url, [Type]::Missing, [Type]::Missing, $url) | Out-Null
    $cell.Style = "Hyperlink"

    $ws.Cells.Item($r,8).Value = $dates[$r - 4]
    $ws.Cells.Item($r,9).Value = $r - 4
}

# --- View-state: selection left on the Contracts sheet ---
$wsContracts = $wb.Worksheets.Item("Contracts")
$wsContracts.Activate() | Out-Null
$wsContracts.Range("B32").Select() | Out-Null

# --- Return focus to Documents (the tab that stays active) and update its view ---
$ws.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 85
$ws.Range("A8").Select() | Out-Null
